$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.091.73"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").Value = "2.495.72"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "'486.00"
$ws.Range("E5").Value = "  +2.66%  "
$ws.Range("D6").Value = "'145.12"
$ws.Range("E6").Value = "  +7.95%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D8").Value = "'0.509"
$ws.Range("E8").Value = "  +2.77%  "
$ws.Range("D9").Value = "2.527.30"
$ws.Range("E9").Value = "  +1.11%  "
$ws.Range("D10").Value = "'5.69"
$ws.Range("E10").Value = "  +3.97%  "
$ws.Range("D11").Value = "'0.0970"
$ws.Range("E11").Value = "  -1.30%  "
$ws.Range("D12").Value = "'0.331"
$ws.Range("E12").Value = "  +2.71%  "
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").Value = "2.953.23"
$ws.Range("E14").Value = "  +1.30%  "
$ws.Range("D15").Value = "56.163.80"
$ws.Range("E15").Value = "  +2.48%  "
$ws.Range("D16").Value = "'21.08"
$ws.Range("E16").Value = "  +4.59%  "
$ws.Range("E17").Value = "  -1.95%  "
$ws.Range("D18").Value = "2.518.93"
$ws.Range("E18").Value = "  +1.55%  "
$ws.Range("D19").Value = "'4.48"
$ws.Range("E19").Value = "  +5.22%  "
$ws.Range("D20").Value = "'10.23"
$ws.Range("E20").Value = "  +6.22%  "
$ws.Range("D21").Value = "'320.18"
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("D23").Value = "'5.78"
$ws.Range("E23").Value = "  +5.92%  "
$ws.Range("D24").Value = "'58.58"
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("D25").Value = "'0.410"
$ws.Range("E25").Value = "  +4.92%  "
$ws.Range("E26").Value = "  +5.17%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("D28").Value = "2.620.21"
$ws.Range("E28").Value = "  +3.32%  "
$ws.Range("D29").Value = "'7.54"
$ws.Range("E29").Value = "  +2.64%  "
$ws.Range("D30").Value = "0.0₃0786"
$ws.Range("E30").Value = "  +3.74%  "
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("D32").Value = "'148.77"
$ws.Range("E32").Value = "  -2.73%  "
$ws.Range("D33").Value = "'18.26"
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("E34").Value = "  +4.60%  "
$ws.Range("D35").Value = "'5.19"
$ws.Range("E35").Value = "  +1.80%  "
$ws.Range("D36").Value = "'1.14"
$ws.Range("E36").Value = "  +5.66%  "
$ws.Range("D37").Value = "'3.69"
$ws.Range("E37").Value = "  +1.06%  "
$ws.Range("D38").Value = "'0.866"
$ws.Range("E38").Value = "  +4.93%  "
$ws.Range("D39").Value = "'34.03"
$ws.Range("E39").Value = "  +0.48%  "
$ws.Range("D40").Value = "'3.53"
$ws.Range("E40").Value = "  +5.72%  "
$ws.Range("D41").Value = "'0.623"
$ws.Range("E41").Value = "  +1.36%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "'0.0556"
$ws.Range("E42").Value = "  +2.83%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'0.994"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("E44").Value = "  +4.03%  "
$ws.Range("D45").Value = "'4.84"
$ws.Range("E45").Value = "  +7.52%  "
$ws.Range("D46").Value = "'260.71"
$ws.Range("E46").Value = "  +14.06%  "
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").Value = "'10.17"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0228"
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("E49").Value = "  +3.24%  "
$ws.Range("D50").Value = "1.927.26"
$ws.Range("E50").Value = "  -2.57%  "
$ws.Range("D51").Value = "'17.58"
$ws.Range("E51").Value = "  +3.83%  "
